$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 11699.385
$ws.Range("J43").Value = 12159.2
$ws.Range("L43").Value = 12159.2
$ws.Range("N43").Value = -12297.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 116212.89
$ws.Range("I62").Value = 146995.86
$ws.Range("K62").Value = 146995.86
$ws.Range("M62").Value = -146371.86

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 116212.89
$ws.Range("I65").Value = 146995.86
$ws.Range("K65").Value = 734979.2999999999
$ws.Range("M65").Value = -731859.2999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4738.5713
$ws.Range("J112").Value = 4738.5713
$ws.Range("L112").Value = 14215.7139
$ws.Range("N112").Value = -16431.7139

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4097.4287
$ws.Range("I132").Value = 1959.1
$ws.Range("J132").Value = 9443.25
$ws.Range("K132").Value = 5877.299999999999
$ws.Range("L132").Value = 28329.75
$ws.Range("M132").Value = -3347.299999999999
$ws.Range("N132").Value = -33389.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1050.381
$ws.Range("I135").Value = 790.5294
$ws.Range("K135").Value = 7114.7646
$ws.Range("M135").Value = -4579.7646

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3823.261
$ws.Range("J137").Value = 10777.363
$ws.Range("L137").Value = 32332.089
$ws.Range("N137").Value = -37432.089

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2371.9355
$ws.Range("I138").Value = 1592.5333
$ws.Range("K138").Value = 4777.5999
$ws.Range("M138").Value = 362.4000999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32771.312
$ws.Range("I32").Value = 18899.305
$ws.Range("J32").Value = 135077.38
$ws.Range("K32").Value = 18899.305
$ws.Range("L32").Value = 135077.38
$ws.Range("M32").Value = -18612.305
$ws.Range("N32").Value = -135651.38

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1162.2632
$ws.Range("I61").Value = 804.4666999999999
$ws.Range("K61").Value = 804.4666999999999
$ws.Range("M61").Value = -592.4666999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2207.7
$ws.Range("I102").Value = 2207.7
$ws.Range("K102").Value = 2207.7
$ws.Range("M102").Value = -585.6999999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1827.8462
$ws.Range("I122").Value = 1703.0526
$ws.Range("J122").Value = 2166.5715
$ws.Range("K122").Value = 5109.1578
$ws.Range("L122").Value = 6499.7145
$ws.Range("M122").Value = -2659.1578
$ws.Range("N122").Value = -11399.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1701.7805
$ws.Range("I132").Value = 1255.069
$ws.Range("J132").Value = 2781.3333
$ws.Range("K132").Value = 3765.207
$ws.Range("L132").Value = 8343.999899999999
$ws.Range("M132").Value = -1235.207
$ws.Range("N132").Value = -13403.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1162.2632
$ws.Range("I136").Value = 804.4666999999999
$ws.Range("K136").Value = 2413.4001
$ws.Range("M136").Value = 136.5999000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 88888
$ws.Range("J137").Value = 88888
$ws.Range("L137").Value = 88888
$ws.Range("N137").Value = -99088

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6951.6772
$ws.Range("I20").Value = 6291.9473
$ws.Range("J20").Value = 7996.25
$ws.Range("K20").Value = 6291.9473
$ws.Range("L20").Value = 7996.25
$ws.Range("M20").Value = -6044.9473
$ws.Range("N20").Value = -8490.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2057.44
$ws.Range("I86").Value = 1400
$ws.Range("J86").Value = 2894.182
$ws.Range("K86").Value = 1400
$ws.Range("L86").Value = 2894.182
$ws.Range("M86").Value = -277
$ws.Range("N86").Value = -5140.182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2057.44
$ws.Range("I89").Value = 1400
$ws.Range("J89").Value = 2894.182
$ws.Range("K89").Value = 7000
$ws.Range("L89").Value = 14470.91
$ws.Range("M89").Value = -1384
$ws.Range("N89").Value = -25702.91

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4854.5654
$ws.Range("I105").Value = 4557.8
$ws.Range("K105").Value = 4557.8
$ws.Range("M105").Value = -2810.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1443.2
$ws.Range("I107").Value = 1554.0625
$ws.Range("K107").Value = 1554.0625
$ws.Range("M107").Value = 365.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1861.8628
$ws.Range("I31").Value = 1273.7222
$ws.Range("J31").Value = 3273.4
$ws.Range("K31").Value = 1273.7222
$ws.Range("L31").Value = 3273.4
$ws.Range("M31").Value = -978.7221999999999
$ws.Range("N31").Value = -3863.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1861.8628
$ws.Range("I34").Value = 1273.7222
$ws.Range("J34").Value = 3273.4
$ws.Range("K34").Value = 1273.7222
$ws.Range("L34").Value = 3273.4
$ws.Range("M34").Value = -1071.7222
$ws.Range("N34").Value = -3677.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 8311.375
$ws.Range("I122").Value = 7999.75
$ws.Range("J122").Value = 8623
$ws.Range("K122").Value = 23999.25
$ws.Range("L122").Value = 25869
$ws.Range("M122").Value = -21549.25
$ws.Range("N122").Value = -30769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2590.7222
$ws.Range("I132").Value = 2289.625
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 6868.875
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -4338.875
$ws.Range("N132").Value = -20058.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1240.0625
$ws.Range("I34").Value = 148.125
$ws.Range("J34").Value = 2332
$ws.Range("K34").Value = 444.375
$ws.Range("L34").Value = 6996
$ws.Range("M34").Value = -360.375
$ws.Range("N34").Value = -7164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3777.8333
$ws.Range("J39").Value = 4792.5
$ws.Range("L39").Value = 14377.5
$ws.Range("N39").Value = -14965.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 9617420
$ws.Range("I55").Value = 591
$ws.Range("J55").Value = 13891566
$ws.Range("K55").Value = 1773
$ws.Range("L55").Value = 41674698
$ws.Range("M55").Value = -1596
$ws.Range("N55").Value = -41675052

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4811.125
$ws.Range("J70").Value = 7166.6665
$ws.Range("L70").Value = 21499.9995
$ws.Range("N70").Value = -22129.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 4811.125
$ws.Range("J73").Value = 7166.6665
$ws.Range("L73").Value = 21499.9995
$ws.Range("N73").Value = -23683.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 9000
$ws.Range("I88").Value = 3000
$ws.Range("K88").Value = 9000
$ws.Range("M88").Value = -8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 9000
$ws.Range("I91").Value = 3000
$ws.Range("K91").Value = 9000
$ws.Range("M91").Value = -7518

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1105.5714
$ws.Range("I122").Value = 907.8
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 8170.2
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -5720.2
$ws.Range("N122").Value = -19300

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 274.6087
$ws.Range("I2").Value = 298.75
$ws.Range("K2").Value = 298.75
$ws.Range("M2").Value = -185.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4047.6
$ws.Range("I70").Value = 3997
$ws.Range("K70").Value = 3997
$ws.Range("M70").Value = -3727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4047.6
$ws.Range("I73").Value = 3997
$ws.Range("K73").Value = 3997
$ws.Range("M73").Value = -3061

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2828.9167
$ws.Range("I126").Value = 2776.5334
$ws.Range("J126").Value = 2916.2222
$ws.Range("K126").Value = 8329.600199999999
$ws.Range("L126").Value = 8748.6666
$ws.Range("M126").Value = -5859.600199999999
$ws.Range("N126").Value = -13688.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 573.1177
$ws.Range("I16").Value = 596.4375
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 596.4375
$ws.Range("L16").Value = 200
$ws.Range("M16").Value = -426.4375
$ws.Range("N16").Value = -540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3848.05
$ws.Range("I40").Value = 2212.7856
$ws.Range("K40").Value = 2212.7856
$ws.Range("M40").Value = -2076.7856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2013.3125
$ws.Range("I46").Value = 2149.125
$ws.Range("J46").Value = 1877.5
$ws.Range("K46").Value = 2149.125
$ws.Range("L46").Value = 1877.5
$ws.Range("M46").Value = -1961.125
$ws.Range("N46").Value = -2253.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 12581.546
$ws.Range("I122").Value = 14344.111
$ws.Range("K122").Value = 43032.333
$ws.Range("M122").Value = -40582.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3757.75
$ws.Range("I136").Value = 2649.75
$ws.Range("K136").Value = 7949.25
$ws.Range("M136").Value = -5399.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4933.222
$ws.Range("J96").Value = 3599.3333
$ws.Range("L96").Value = 3599.3333
$ws.Range("N96").Value = -6345.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 60000
$ws.Range("J108").Value = 60000
$ws.Range("L108").Value = 60000
$ws.Range("N108").Value = -67680

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2132.45
$ws.Range("I126").Value = 1795.8
$ws.Range("J126").Value = 3142.4
$ws.Range("K126").Value = 5387.4
$ws.Range("L126").Value = 9427.200000000001
$ws.Range("M126").Value = -2917.4
$ws.Range("N126").Value = -14367.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 21788.182
$ws.Range("I132").Value = 21240.055
$ws.Range("J132").Value = 24254.75
$ws.Range("K132").Value = 63720.165
$ws.Range("L132").Value = 72764.25
$ws.Range("M132").Value = -61190.165
$ws.Range("N132").Value = -77824.25
